$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 333584.84
$ws.Range("I12").Value = 236.33333
$ws.Range("J12").Value = 666933.3
$ws.Range("K12").Value = 236.33333
$ws.Range("L12").Value = 666933.3
$ws.Range("M12").Value = -66.33332999999999
$ws.Range("N12").Value = -667273.3
$ws.Range("H69").Value = 3108.6667
$ws.Range("I69").Value = 3108.6667
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 9326.000100000001
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -8452.000100000001
$ws.Range("H72").Value = 3108.6667
$ws.Range("I72").Value = 3108.6667
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 27978.0003
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -23610.0003
$ws.Range("H96").Value = 870.6923
$ws.Range("I96").Value = 718.3333
$ws.Range("J96").Value = 1001.2857
$ws.Range("K96").Value = 2154.9999
$ws.Range("L96").Value = 3003.8571
$ws.Range("M96").Value = -781.9998999999998
$ws.Range("N96").Value = -5749.8571
$ws.Range("H129").Value = 2632707
$ws.Range("I129").Value = 125003000
$ws.Range("J129").Value = 1087.914
$ws.Range("K129").Value = 375009000
$ws.Range("L129").Value = 3263.742
$ws.Range("M129").Value = -375004000
$ws.Range("N129").Value = -13263.742
$ws.Range("H135").Value = 846.8333
$ws.Range("I135").Value = 624.09753
$ws.Range("J135").Value = 2151.4285
$ws.Range("K135").Value = 5616.87777
$ws.Range("L135").Value = 19362.8565
$ws.Range("M135").Value = -3081.87777
$ws.Range("N135").Value = -24432.8565
$ws.Range("H138").Value = 3277.4443
$ws.Range("I138").Value = 1402.1613
$ws.Range("J138").Value = 5805
$ws.Range("K138").Value = 4206.4839
$ws.Range("L138").Value = 17415
$ws.Range("M138").Value = 933.5160999999998
$ws.Range("N138").Value = -27695
$ws.Range("H141").Value = 437535.84
$ws.Range("I141").Value = 2965.182
$ws.Range("J141").Value = 6412882
$ws.Range("K141").Value = 8895.545999999998
$ws.Range("L141").Value = 19238646
$ws.Range("M141").Value = -3715.545999999998
$ws.Range("N141").Value = -19249006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2964.3333
$ws.Range("I61").Value = 1128.6666
$ws.Range("J61").Value = 4800
$ws.Range("K61").Value = 1128.6666
$ws.Range("L61").Value = 4800
$ws.Range("M61").Value = -916.6666
$ws.Range("N61").Value = -5224
$ws.Range("H63").Value = 3272.182
$ws.Range("I63").Value = 1883.6923
$ws.Range("J63").Value = 5277.778
$ws.Range("K63").Value = 1883.6923
$ws.Range("L63").Value = 5277.778
$ws.Range("M63").Value = -1197.6923
$ws.Range("N63").Value = -6649.778
$ws.Range("H66").Value = 3272.182
$ws.Range("I66").Value = 1883.6923
$ws.Range("J66").Value = 5277.778
$ws.Range("K66").Value = 9418.461499999999
$ws.Range("L66").Value = 26388.89
$ws.Range("M66").Value = -5986.461499999999
$ws.Range("N66").Value = -33252.89
$ws.Range("H74").Value = 2000
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 1000
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -126
$ws.Range("N74").Value = -4748
$ws.Range("H77").Value = 2000
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 5000
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -632
$ws.Range("N77").Value = -23736
$ws.Range("H122").Value = 5166.6665
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 7333.3335
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 22000.0005
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -26900.0005
$ws.Range("H132").Value = 3183.9285
$ws.Range("I132").Value = 2173.0557
$ws.Range("J132").Value = 5003.5
$ws.Range("K132").Value = 6519.1671
$ws.Range("L132").Value = 15010.5
$ws.Range("M132").Value = -3989.1671
$ws.Range("N132").Value = -20070.5
$ws.Range("H136").Value = 2964.3333
$ws.Range("I136").Value = 1128.6666
$ws.Range("J136").Value = 4800
$ws.Range("K136").Value = 3385.9998
$ws.Range("L136").Value = 14400
$ws.Range("M136").Value = -835.9998000000001
$ws.Range("N136").Value = -19500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 29986.666
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 29986.666
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 29986.666
$ws.Range("N130").Value = -40026.666
$ws.Range("H134").Value = 6231.8237
$ws.Range("I134").Value = 10537.167
$ws.Range("J134").Value = 3883.4546
$ws.Range("K134").Value = 31611.501
$ws.Range("L134").Value = 11650.3638
$ws.Range("M134").Value = -29076.501
$ws.Range("N134").Value = -16720.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 31.416666
$ws.Range("I7").Value = 12.428572
$ws.Range("J7").Value = 58
$ws.Range("K7").Value = 12.428572
$ws.Range("L7").Value = 58
$ws.Range("M7").Value = 100.571428
$ws.Range("N7").Value = -284
$ws.Range("H31").Value = 2711.4038
$ws.Range("I31").Value = 1635.7179
$ws.Range("J31").Value = 5938.4614
$ws.Range("K31").Value = 1635.7179
$ws.Range("L31").Value = 5938.4614
$ws.Range("M31").Value = -1340.7179
$ws.Range("N31").Value = -6528.4614
$ws.Range("H34").Value = 2711.4038
$ws.Range("I34").Value = 1635.7179
$ws.Range("J34").Value = 5938.4614
$ws.Range("K34").Value = 1635.7179
$ws.Range("L34").Value = 5938.4614
$ws.Range("M34").Value = -1433.7179
$ws.Range("N34").Value = -6342.4614
$ws.Range("H74").Value = 20986.7
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 20986.7
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 20986.7
$ws.Range("N74").Value = -22734.7
$ws.Range("H77").Value = 20986.7
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 20986.7
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 62960.10000000001
$ws.Range("N77").Value = -71696.10000000001
$ws.Range("H94").Value = 1367.8182
$ws.Range("I94").Value = 1432.3334
$ws.Range("J94").Value = 1353.4814
$ws.Range("K94").Value = 1432.3334
$ws.Range("L94").Value = 1353.4814
$ws.Range("M94").Value = -981.3334
$ws.Range("N94").Value = -2255.4814
$ws.Range("H122").Value = 4010.85
$ws.Range("I122").Value = 3018.1667
$ws.Range("J122").Value = 5499.875
$ws.Range("K122").Value = 9054.500100000001
$ws.Range("L122").Value = 16499.625
$ws.Range("M122").Value = -6604.500100000001
$ws.Range("N122").Value = -21399.625
$ws.Range("H132").Value = 2050.468
$ws.Range("I132").Value = 1652.8125
$ws.Range("J132").Value = 2898.8
$ws.Range("K132").Value = 4958.4375
$ws.Range("L132").Value = 8696.400000000001
$ws.Range("M132").Value = -2428.4375
$ws.Range("N132").Value = -13756.4
$ws.Range("H134").Value = 1453.6451
$ws.Range("I134").Value = 996.4897999999999
$ws.Range("J134").Value = 3176.7693
$ws.Range("K134").Value = 2989.4694
$ws.Range("L134").Value = 9530.3079
$ws.Range("M134").Value = -454.4694
$ws.Range("N134").Value = -14600.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1314.5
$ws.Range("I5").Value = 636.6
$ws.Range("J5").Value = 1691.1111
$ws.Range("K5").Value = 1909.8
$ws.Range("L5").Value = 5073.3333
$ws.Range("M5").Value = -1797.8
$ws.Range("N5").Value = -5297.3333
$ws.Range("H131").Value = 2088.889
$ws.Range("I131").Value = 3128.75
$ws.Range("J131").Value = 1651.0526
$ws.Range("K131").Value = 9386.25
$ws.Range("L131").Value = 4953.1578
$ws.Range("M131").Value = -4346.25
$ws.Range("N131").Value = -15033.1578
$ws.Range("H132").Value = 2161.077
$ws.Range("I132").Value = 1348.625
$ws.Range("J132").Value = 3461
$ws.Range("K132").Value = 12137.625
$ws.Range("L132").Value = 31149
$ws.Range("M132").Value = -9607.625
$ws.Range("N132").Value = -36209
$ws.Range("H135").Value = 1314.5
$ws.Range("I135").Value = 636.6
$ws.Range("J135").Value = 1691.1111
$ws.Range("K135").Value = 5729.400000000001
$ws.Range("L135").Value = 15219.9999
$ws.Range("M135").Value = -3194.400000000001
$ws.Range("N135").Value = -20289.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3685.1738
$ws.Range("I122").Value = 1996.3636
$ws.Range("J122").Value = 5233.25
$ws.Range("K122").Value = 5989.0908
$ws.Range("L122").Value = 15699.75
$ws.Range("M122").Value = -3539.0908
$ws.Range("N122").Value = -20599.75
$ws.Range("H126").Value = 3776.2222
$ws.Range("I126").Value = 1996.5
$ws.Range("J126").Value = 5200
$ws.Range("K126").Value = 5989.5
$ws.Range("L126").Value = 15600
$ws.Range("M126").Value = -3519.5
$ws.Range("N126").Value = -20540
$ws.Range("H132").Value = 3195.0417
$ws.Range("I132").Value = 2975.7144
$ws.Range("J132").Value = 3785.5386
$ws.Range("K132").Value = 8927.143199999999
$ws.Range("L132").Value = 11356.6158
$ws.Range("M132").Value = -6397.143199999999
$ws.Range("N132").Value = -16416.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 76923930
$ws.Range("I16").Value = 111112220
$ws.Range("J16").Value = 268
$ws.Range("K16").Value = 111112220
$ws.Range("L16").Value = 268
$ws.Range("M16").Value = -111112050
$ws.Range("N16").Value = -608
$ws.Range("H132").Value = 2597.8276
$ws.Range("I132").Value = 1185.5
$ws.Range("J132").Value = 4908.909
$ws.Range("K132").Value = 3556.5
$ws.Range("L132").Value = 14726.727
$ws.Range("M132").Value = -1026.5
$ws.Range("N132").Value = -19786.727
$ws.Range("H136").Value = 2141.8918
$ws.Range("I136").Value = 1522.4
$ws.Range("J136").Value = 4796.857
$ws.Range("K136").Value = 4567.200000000001
$ws.Range("L136").Value = 14390.571
$ws.Range("M136").Value = -2017.200000000001
$ws.Range("N136").Value = -19490.571
$ws.Range("H139").Value = 49736.25
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 49736.25
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 49736.25
$ws.Range("N139").Value = -60016.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1119.1111
$ws.Range("I96").Value = 1013
$ws.Range("J96").Value = 1331.3334
$ws.Range("K96").Value = 1013
$ws.Range("L96").Value = 1331.3334
$ws.Range("M96").Value = 360
$ws.Range("N96").Value = -4077.3334
$ws.Range("H122").Value = 529085.5
$ws.Range("I122").Value = 835443.3
$ws.Range("J122").Value = 3900.5715
$ws.Range("K122").Value = 2506329.9
$ws.Range("L122").Value = 11701.7145
$ws.Range("M122").Value = -2503879.9
$ws.Range("N122").Value = -16601.7145
$ws.Range("H132").Value = 7612.07
$ws.Range("I132").Value = 1788.9246
$ws.Range("J132").Value = 16964.395
$ws.Range("K132").Value = 5366.7738
$ws.Range("L132").Value = 50893.185
$ws.Range("M132").Value = -2836.7738
$ws.Range("N132").Value = -55953.185
